# Applies the weekly odds-update changes to Sheet1, matching the
# FlashScore workbook diff (values for 2025-02-24 fixtures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.04
$ws.Range("O2").Value = 1.22
$ws.Range("S2").Value = 2.75
$ws.Range("T2").Value = 1.44

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("M3").Value = 1.11
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.63
$ws.Range("T3").Value = 1.17
$ws.Range("W3").Value = 2.2
$ws.Range("X3").Value = 1.62

# Row 4
$ws.Range("M4").Value = 1.13
$ws.Range("O4").Value = 1.57
$ws.Range("T4").Value = 1.13

# Row 5
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 5.25
$ws.Range("K5").Value = 1.95
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("O5").Value = 1.5
$ws.Range("T5").Value = 1.17
$ws.Range("AA5").Value = 9.5
$ws.Range("AE5").Value = 6.5
$ws.Range("AO5").Value = 67

# Row 6
$ws.Range("M6").Value = 1.13
$ws.Range("O6").Value = 1.57
$ws.Range("T6").Value = 1.13

# Row 7
$ws.Range("G7").Value = 2.2
$ws.Range("I7").Value = 3.9
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 5
$ws.Range("AL7").Value = 15

# Row 8
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 6

# Row 9
$ws.Range("H9").Value = 3.1
$ws.Range("J9").Value = 3.2
$ws.Range("K9").Value = 2.05
$ws.Range("N9").Value = 8.5
$ws.Range("U9").Value = 1.44
$ws.Range("V9").Value = 2.63
$ws.Range("W9").Value = 1.83
$ws.Range("X9").Value = 1.83
$ws.Range("Y9").Value = 7.5
$ws.Range("AE9").Value = 8.5
$ws.Range("AG9").Value = 15
$ws.Range("AI9").Value = 301
$ws.Range("AJ9").Value = 8.5
$ws.Range("AN9").Value = 23
$ws.Range("AO9").Value = 34

# Row 11
$ws.Range("G11").Value = 1.36
$ws.Range("L11").Value = 8.5
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.5
$ws.Range("Y11").Value = 6
$ws.Range("AB11").Value = 8
$ws.Range("AJ11").Value = 17
$ws.Range("AP11").Value = 1.47
$ws.Range("AQ11").Value = 2.65
$ws.Range("AR11").Value = 2.65
$ws.Range("AS11").Value = 1.47

# Row 15
$ws.Range("G15").Value = 2.15
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 2.75
$ws.Range("W15").Value = 1.62
$ws.Range("X15").Value = 2.2
$ws.Range("AG15").Value = 13

# Row 16
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("Q16").Value = 2.08
$ws.Range("R16").Value = 1.73

# Row 17
$ws.Range("G17").Value = 3.7
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 1.91
$ws.Range("AK17").Value = 9
$ws.Range("AM17").Value = 17
$ws.Range("AN17").Value = 17

# Row 18
$ws.Range("H18").Value = 3.75
$ws.Range("P18").Value = 3.2
$ws.Range("W18").Value = 1.8
$ws.Range("X18").Value = 1.8
$ws.Range("Y18").Value = 6.7
$ws.Range("Z18").Value = 7.3
$ws.Range("AB18").Value = 11.25
$ws.Range("AD18").Value = 26
$ws.Range("AF18").Value = 7.4
$ws.Range("AJ18").Value = 14.5
$ws.Range("AN18").Value = 55
$ws.Range("AO18").Value = 55

# Row 20
$ws.Range("G20").Value = 1.9

# Row 21
$ws.Range("G21").Value = 1.02
$ws.Range("H21").Value = 12
$ws.Range("I21").Value = 23
$ws.Range("J21").Value = 1.14
$ws.Range("K21").Value = 5.3
$ws.Range("L21").Value = 16
$ws.Range("Y21").Value = 65
$ws.Range("Z21").Value = 26
$ws.Range("AA21").Value = 32
$ws.Range("AB21").Value = 14.5
$ws.Range("AC21").Value = 18
$ws.Range("AD21").Value = 40
$ws.Range("AE21").Value = 150
$ws.Range("AF21").Value = 80
$ws.Range("AG21").Value = 75
$ws.Range("AH21").Value = 120
$ws.Range("AI21").Value = 400
$ws.Range("AJ21").Value = 450
$ws.Range("AL21").Value = 200
$ws.Range("AO21").Value = 250

# Row 22
$ws.Range("G22").Value = 2.7
$ws.Range("I22").Value = 2.4

# Row 23
$ws.Range("I23").Value = 1.85

# Row 24
$ws.Range("I24").Value = 1.75

# Row 25
$ws.Range("G25").Value = 2.38

# Row 26
$ws.Range("G26").Value = 2

# Row 27
$ws.Range("I27").Value = 2.88

# Row 29
$ws.Range("G29").Value = 1.67
$ws.Range("M29").Value = 1.07
$ws.Range("N29").Value = 9
$ws.Range("Q29").Value = 2.15
$ws.Range("R29").Value = 1.67
$ws.Range("S29").Value = 4
$ws.Range("T29").Value = 1.22

# Row 30
$ws.Range("G30").Value = 2.3
$ws.Range("I30").Value = 2.88
$ws.Range("M30").Value = 1.02
$ws.Range("O30").Value = 1.25
$ws.Range("T30").Value = 1.36
$ws.Range("Y30").Value = 9
$ws.Range("AA30").Value = 9.5
$ws.Range("AB30").Value = 21
$ws.Range("AJ30").Value = 11
$ws.Range("AN30").Value = 23

# Row 31
$ws.Range("M31").Value = 1.03
$ws.Range("O31").Value = 1.25
$ws.Range("T31").Value = 1.33

# Row 33
$ws.Range("M33").Value = 1.02
$ws.Range("O33").Value = 1.13
$ws.Range("T33").Value = 1.73
